$d = $word.ActiveDocument

$replacements = @(
    @{old = "774×9=6966"; new = "113×8=904"},
    @{old = "518×3=1554"; new = "429×9=3861"},
    @{old = "688×3=2064"; new = "290×8=2320"},
    @{old = "588×4=2352"; new = "263×4=1052"},
    @{old = "810×5=4050"; new = "111×6=666"},
    @{old = "932×3=2796"; new = "603×7=4221"},
    @{old = "334×4=1336"; new = "707×3=2121"},
    @{old = "662×3=1986"; new = "107×2=214"},
    @{old = "542×3=1626"; new = "786×4=3144"},
    @{old = "150×2=300"; new = "996×2=1992"},
    @{old = "437×6=2622"; new = "218×5=1090"},
    @{old = "799×4=3196"; new = "212×9=1908"},
    @{old = "854×4=3416"; new = "323×4=1292"},
    @{old = "819×9=7371"; new = "830×5=4150"},
    @{old = "317×5=1585"; new = "920×3=2760"},
    @{old = "917×2=1834"; new = "831×5=4155"},
    @{old = "521×7=3647"; new = "848×9=7632"},
    @{old = "825×5=4125"; new = "467×2=934"},
    @{old = "792×9=7128"; new = "892×4=3568"},
    @{old = "745×8=5960"; new = "946×9=8514"},
    @{old = "644×5=3220"; new = "358×6=2148"},
    @{old = "659×7=4613"; new = "644×4=2576"},
    @{old = "964×2=1928"; new = "133×2=266"},
    @{old = "437×7=3059"; new = "584×7=4088"},
    @{old = "984×8=7872"; new = "741×2=1482"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $r.new, 2)
}
